$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# The cell E8 held the "Good Morning" greeting; update it per the
# outside-Webstudio git commit ("GIT UPDATE").
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the resulting cursor/selection position (matches the saved
# sheetView's <selection activeCell="E8" sqref="E8"/>).
$ws.Range("E8").Select()
